$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values look numeric need to stay text (matching
# the original inlineStr cell type), so force a text number format before assignment.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "69.224.90"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "3.752.80"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "602.36"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").Value = "167.12"
$ws.Range("E6").Value = "  -1.41%  "
$ws.Range("D7").Value = "3.752.21"
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "0.539"
$ws.Range("E9").Value = "  +1.20%  "
$ws.Range("E10").Value = "  +2.25%  "
$ws.Range("D11").Value = "6.40"
$ws.Range("E11").Value = "  +1.37%  "
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("D13").Value = "38.08"
$ws.Range("E13").Value = "  -1.10%  "
$ws.Range("E14").Value = "  +1.49%  "
$ws.Range("D15").Value = "4.381.80"
$ws.Range("E15").Value = "  +0.32%  "
$ws.Range("D16").Value = "3.746.01"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").Value = "69.216.13"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").Value = "7.39"
$ws.Range("E18").Value = "  +1.40%  "
$ws.Range("D19").Value = "17.43"
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("D21").Value = "11.18"
$ws.Range("E21").Value = "  +18.54%  "
$ws.Range("D22").Value = "494.42"
$ws.Range("E22").Value = "  -0.87%  "
$ws.Range("D23").Value = "0.728"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("E24").Value = "  +7.40%  "
$ws.Range("D25").Value = "84.91"
$ws.Range("E25").Value = "  -0.18%  "
$ws.Range("E26").Value = "  -0.79%  "
$ws.Range("D27").Value = "12.32"
$ws.Range("E27").Value = "  -0.26%  "
$ws.Range("D28").Value = "10.11"
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  +1.81%  "
$ws.Range("E31").Value = "  +1.84%  "
$ws.Range("D32").Value = "8.08"
$ws.Range("E32").Value = "  +0.46%  "
$ws.Range("D33").Value = "31.63"
$ws.Range("E33").Value = "  -0.42%  "
$ws.Range("D34").Value = "3.896.90"
$ws.Range("E34").Value = "  +0.36%  "
$ws.Range("D35").Value = "3.691.43"
$ws.Range("E35").Value = "  +0.24%  "
$ws.Range("D36").Value = "0.109"
$ws.Range("E36").Value = "  -0.54%  "
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D38").Value = "5.99"
$ws.Range("E38").Value = "  +2.94%  "
$ws.Range("B39").Value = "Mantle"
$ws.Range("C39").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D39").Value = "1.02"
$ws.Range("E39").Value = "  +0.94%  "
$ws.Range("E40").Value = "  +0.83%  "
$ws.Range("D41").Value = "0.325"
$ws.Range("E41").Value = "  -0.38%  "
$ws.Range("D42").Value = "2.99"
$ws.Range("E42").Value = "  +4.48%  "
$ws.Range("D43").Value = "430.77"
$ws.Range("E43").Value = "  -1.66%  "
$ws.Range("D44").Value = "48.82"
$ws.Range("E44").Value = "  -0.73%  "
$ws.Range("D46").Value = "8.49"
$ws.Range("E46").Value = "  +0.91%  "
$ws.Range("D48").Value = "40.26"
$ws.Range("E48").Value = "  -0.75%  "
$ws.Range("D49").Value = "141.29"
$ws.Range("E49").Value = "  -1.03%  "
$ws.Range("D50").Value = "2.793.03"
$ws.Range("E50").Value = "  +1.16%  "
$ws.Range("D51").Value = "0.0353"
$ws.Range("E51").Value = "  -0.06%  "
